$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header text in L1 (typo change: OBSERVACIÓN -> OBSEVACIÓN)
$ws.Range("L1").Value = "OBSEVACIÓN"

# Update the active selection to G1 as in the edited file
$ws.Range("G1").Select()
